$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace resource data (LGBTQ Trevor Project -> NYS Problem Gambling) ---
$ws.Range("A2").Value = "GAMB_NYSmainwebsite"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "NYS Problem Gambling Help"
$ws.Range("C2").Value = "GAMBdata()"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "NYS Problem Gambling Help"
$ws.Range("E2").Value = "F"
$ws.Range("F2").Value = "T"

# Pull the freshly recalculated generated-code strings (B4/B5/B6 formulas depend on row 2)
$b4 = $ws.Range("B4").Value2
$b5 = $ws.Range("B5").Value2
$b6 = $ws.Range("B6").Value2

# --- Rows 10-15: drop old leftover reference/staging data, shift generated code up ---
$ws.Range("A10").Value = $b4
$ws.Range("C10").Clear()
$ws.Range("D10").Clear()
$ws.Range("E10").Clear()
$ws.Range("H10").Clear()
$ws.Range("G10").ClearContents()
$ws.Range("I10").ClearContents()

$ws.Range("A11").Value = $b5

$ws.Range("A12").Value = $b6

$ws.Range("A13").ClearContents()
$ws.Range("A14").Clear()
$ws.Range("A15").Clear()

# --- view state (best effort) ---
$ws.Range("D5").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$win.Left = 27860
$win.Top = 8720
